$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 261, shifting existing rows 261:273 down to 262:274
$ws.Rows.Item(261).Insert()

# Populate the newly inserted row 261 with the new record (copy of the
# previous row 261 template, with updated Fecha/Volumen/Precio values)
$ws.Cells.Item(261, 1).Value = 11
$ws.Cells.Item(261, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(261, 3).Value = "Bíobío"
$ws.Cells.Item(261, 4).Value = 44747
$ws.Cells.Item(261, 5).Value = 8
$ws.Cells.Item(261, 6).Value = 100112009
$ws.Cells.Item(261, 7).Value = "Acelga"
$ws.Cells.Item(261, 8).Value = "Sin especificar"
$ws.Cells.Item(261, 9).Value = "Primera"
$ws.Cells.Item(261, 10).Value = 450
$ws.Cells.Item(261, 11).Value = 600
$ws.Cells.Item(261, 12).Value = 650
$ws.Cells.Item(261, 13).Value = 622
$ws.Cells.Item(261, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(261, 15).Value = "Región de Ñuble"
$ws.Cells.Item(261, 16).Value = 622
$ws.Cells.Item(261, 17).Value = 1
$ws.Cells.Item(261, 18).Value = "Hortaliza"

# Apply the same number format (date format) as the other D column cells
$ws.Cells.Item(261, 4).NumberFormat = $ws.Cells.Item(262, 4).NumberFormat
